{"js": "// Update the two-digit-by-two-digit multiplication problems in the table.\n// Each \"old\" problem text is unique in the document, so a scoped search +\n// in-place text replace keeps every other paragraph/run property untouched.\nconst replacements = [\n  [\"18\u00d711=\", \"22\u00d761=\"],\n  [\"13\u00d714=\", \"71\u00d775=\"],\n  [\"48\u00d791=\", \"76\u00d799=\"],\n  [\"67\u00d711=\", \"74\u00d774=\"],\n  [\"79\u00d767=\", \"68\u00d777=\"],\n  [\"24\u00d726=\", \"56\u00d718=\"],\n  [\"47\u00d756=\", \"84\u00d726=\"],\n  [\"54\u00d755=\", \"65\u00d745=\"],\n  [\"16\u00d786=\", \"78\u00d717=\"],\n  [\"48\u00d742=\", \"98\u00d781=\"],\n  [\"44\u00d777=\", \"77\u00d792=\"],\n  [\"76\u00d763=\", \"14\u00d721=\"],\n  [\"93\u00d756=\", \"82\u00d757=\"],\n  [\"86\u00d733=\", \"82\u00d789=\"],\n  [\"99\u00d781=\", \"32\u00d741=\"],\n  [\"33\u00d723=\", \"96\u00d713=\"],\n  [\"77\u00d790=\", \"61\u00d759=\"],\n  [\"94\u00d731=\", \"83\u00d781=\"],\n  [\"86\u00d777=\", \"26\u00d767=\"],\n  [\"84\u00d737=\", \"87\u00d726=\"],\n  [\"87\u00d763=\", \"85\u00d726=\"],\n  [\"96\u00d761=\", \"31\u00d751=\"],\n  [\"79\u00d787=\", \"36\u00d791=\"],\n  [\"76\u00d791=\", \"81\u00d770=\"],\n  [\"60\u00d753=\", \"45\u00d771=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit-by-two-digit multiplication problems in the table.\n# Each \"old\" problem text is unique in the document, so Find/Replace on the\n# whole document content swaps just the digits and leaves every other run\n# property (font, size, justification, etc.) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"18\u00d711=\", \"22\u00d761=\"),\n    @(\"13\u00d714=\", \"71\u00d775=\"),\n    @(\"48\u00d791=\", \"76\u00d799=\"),\n    @(\"67\u00d711=\", \"74\u00d774=\"),\n    @(\"79\u00d767=\", \"68\u00d777=\"),\n    @(\"24\u00d726=\", \"56\u00d718=\"),\n    @(\"47\u00d756=\", \"84\u00d726=\"),\n    @(\"54\u00d755=\", \"65\u00d745=\"),\n    @(\"16\u00d786=\", \"78\u00d717=\"),\n    @(\"48\u00d742=\", \"98\u00d781=\"),\n    @(\"44\u00d777=\", \"77\u00d792=\"),\n    @(\"76\u00d763=\", \"14\u00d721=\"),\n    @(\"93\u00d756=\", \"82\u00d757=\"),\n    @(\"86\u00d733=\", \"82\u00d789=\"),\n    @(\"99\u00d781=\", \"32\u00d741=\"),\n    @(\"33\u00d723=\", \"96\u00d713=\"),\n    @(\"77\u00d790=\", \"61\u00d759=\"),\n    @(\"94\u00d731=\", \"83\u00d781=\"),\n    @(\"86\u00d777=\", \"26\u00d767=\"),\n    @(\"84\u00d737=\", \"87\u00d726=\"),\n    @(\"87\u00d763=\", \"85\u00d726=\"),\n    @(\"96\u00d761=\", \"31\u00d751=\"),\n    @(\"79\u00d787=\", \"36\u00d791=\"),\n    @(\"76\u00d791=\", \"81\u00d770=\"),\n    @(\"60\u00d753=\", \"45\u00d771=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\nWrite-Output \"done\"\n"}
